$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "maa://24702 (94.2), maa://25390 (95.85), maa://36681 (86.3)"
$ws.Range("L2").Value = "*maa://24633 (55.26), *maa://30515 (69.0), *maa://34787 (71.88), ***maa://20792 (11.93), maa://39402 (84.38), ***maa://29083 (27.78)"
$ws.Range("T2").Value = "maa://22742 (91.67), *maa://20791 (62.32)"
$ws.Range("T3").Value = "maa://24617 (88.57), **maa://20790 (43.94), ***maa://37170 (20.0)"
$ws.Range("X4").Value = "**maa://32495 (47.27), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (87.5)"
$ws.Range("A8").Value = "更新日期：2024.11.22 13:18:20"
$ws.Range("AF8").Value = "*maa://24479 (78.21), *maa://21990 (53.85)"
$ws.Range("X9").Value = "maa://26223 (97.37)"
$ws.Range("D11").Value = "maa://36707 (99.38)"
$ws.Range("AB12").Value = "maa://23669 (95.27), maa://36677 (93.48), maa://39872 (90.0)"
$ws.Range("L14").Value = "maa://26245 (96.3), maa://21288 (96.21), maa://36682 (97.3), maa://39841 (94.12)"
$ws.Range("D15").Value = "*maa://22743 (77.13), maa://22734 (83.76), *maa://30808 (63.93), ***maa://36048 (26.83)"
$ws.Range("AB16").Value = "maa://26228 (95.35)"
$ws.Range("X18").Value = "maa://21917 (97.7), maa://22741 (83.33)"
$ws.Range("H20").Value = "maa://22864 (88.65)"
$ws.Range("X21").Value = "maa://20110 (86.76), maa://34946 (92.31)"
$ws.Range("AF22").Value = "maa://29658 (93.18)"
$ws.Range("L23").Value = "maa://39756 (93.06), maa://39875 (93.22)"
$ws.Range("T23").Value = "maa://24387 (81.08), maa://31212 (96.3)"
$ws.Range("AB23").Value = "maa://29652 (97.5)"
$ws.Range("X25").Value = "*maa://29890 (75.61)"
$ws.Range("AB25").Value = "maa://31215 (84.95), *maa://24516 (79.07), maa://26001 (87.27)"
$ws.Range("X26").Value = "maa://24389 (96.15)"
$ws.Range("AF27").Value = "maa://24023 (97.01)"
$ws.Range("X28").Value = "maa://39929 (89.27), ***maa://39723 (14.29), maa://41749 (85.71)"
$ws.Range("AF28").Value = "maa://36660 (92.49), *maa://36701 (62.96)"
$ws.Range("T30").Value = "*maa://32940 (66.67), maa://24388 (94.44)"
$ws.Range("T32").Value = "maa://41108 (87.5), maa://42859 (93.62), maa://41238 (95.0)"
$ws.Range("T35").Value = "maa://24842 (94.0)"
$ws.Range("T36").Value = "maa://27613 (99.02)"
$ws.Range("T38").Value = "maa://30713 (96.67)"
$ws.Range("AF38").Value = "maa://36697 (85.71)"
$ws.Range("P39").Value = "maa://24709 (91.45)"
$ws.Range("H47").Value = "maa://27410 (96.01), maa://29661 (97.78), maa://28038 (84.62)"
